$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update wall-thickness t[T] parameter (row 13)
$ws.Range("C13").Value = 1

# Update sweep parameters N[0], N[1], dN (rows 36-38)
$ws.Range("C36").Value = 0
$ws.Range("C37").Value = 20
$ws.Range("C38").Value = 1

# Move the active selection to D38
$ws.Range("D38").Select()
